# Update Global_M2 / Hungary sheet with refreshed TradingView data feed values
# and append two new monthly observations (2023-03-01 and 2023-04-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing cells (columns B and D) for rows 358-387 -----------
$updates = @(
    @{Row=358; Col="B"; Value=30734900000000},
    @{Row=358; Col="D"; Value=99090818231.3513},
    @{Row=359; Col="B"; Value=31322400000000},
    @{Row=359; Col="D"; Value=99561353320.36465},
    @{Row=361; Col="B"; Value=33495800000000},
    @{Row=361; Col="D"; Value=112815417570.4259},
    @{Row=362; Col="B"; Value=32988400000000},
    @{Row=362; Col="D"; Value=111980718965.3417},
    @{Row=363; Col="B"; Value=33552400000000},
    @{Row=363; Col="D"; Value=111796614687.4583},
    @{Row=364; Col="B"; Value=33970200000000},
    @{Row=364; Col="D"; Value=110071285075.4974},
    @{Row=365; Col="B"; Value=33881500000000},
    @{Row=365; Col="D"; Value=113126878130.217},
    @{Row=366; Col="B"; Value=33795100000000},
    @{Row=366; Col="D"; Value=118994736008.1689},
    @{Row=367; Col="B"; Value=34176400000000},
    @{Row=367; Col="D"; Value=115377006566.1766},
    @{Row=368; Col="B"; Value=34635200000000},
    @{Row=368; Col="D"; Value=114646232270.2372},
    @{Row=369; Col="B"; Value=34976400000000},
    @{Row=369; Col="D"; Value=118284995400.6818},
    @{Row=370; Col="B"; Value=35509700000000},
    @{Row=370; Col="D"; Value=114397595407.3053},
    @{Row=371; Col="B"; Value=36350700000000},
    @{Row=371; Col="D"; Value=116727517942.2956},
    @{Row=372; Col="B"; Value=37653400000000},
    @{Row=372; Col="D"; Value=117063267526.8149},
    @{Row=373; Col="B"; Value=38869800000000},
    @{Row=373; Col="D"; Value=119742338545.7099},
    @{Row=375; Col="B"; Value=39527800000000},
    @{Row=375; Col="D"; Value=119368847013.3478},
    @{Row=376; Col="B"; Value=39600300000000},
    @{Row=376; Col="D"; Value=119319342543.3509},
    @{Row=377; Col="B"; Value=39874100000000},
    @{Row=377; Col="D"; Value=111147317073.1707},
    @{Row=380; Col="B"; Value=41209500000000},
    @{Row=380; Col="D"; Value=104230214735.5642},
    @{Row=381; Col="B"; Value=42040500000000},
    @{Row=381; Col="D"; Value=105465104610.9076},
    @{Row=382; Col="B"; Value=41918300000000},
    @{Row=382; Col="D"; Value=97100532777.39171},
    @{Row=387; Col="B"; Value=40513200000000},
    @{Row=387; Col="D"; Value=113226720625.1432}
)

foreach ($item in $updates) {
    $addr = "$($item.Col)$($item.Row)"
    $ws.Range($addr).Value = $item.Value
}

# --- 2. Append two new rows (388, 389) with the same row layout/style as
#        the existing data rows -------------------------------------------
$newRows = @(
    @{Row=388; A=44986; B=40410200000000;    C=0.00285257873117298;  D=115273277042.4464},
    @{Row=389; A=45017; B=39781400000000;    C=0.002952770436862386; D=117465341856.9973}
)

foreach ($row in $newRows) {
    $r = $row.Row
    $prev = $r - 1

    # Copy formatting (style) from the row above so the date cell keeps the
    # same number format / border / alignment as the rest of column A.
    $ws.Range("A$prev").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
}
